$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the header style (bold / centered / bordered) that already lives
# on D1 so it can be stamped onto the newly-created header + index cells.
$ws.Range("D1").Copy()

# --- Row 1: header row, now A1:F1 -----------------------------------------
# Existing columns B1/C1/D1 keep their text (C1/D1 unchanged); B1 becomes
# "datname" and a brand new "datnum" column is inserted at A1. Two more
# headers (x_label / y_label) are appended at E1 / F1.
$ws.Range("A1").Value = "datnum"
$ws.Range("B1").Value = "datname"
$ws.Range("C1").Value = "time"
$ws.Range("D1").Value = "picklepath"
$ws.Range("E1").Value = "x_label"
$ws.Range("F1").Value = "y_label"

$ws.Range("A1:F1").PasteSpecial(-4122)

# --- Row 2: first data row --------------------------------------------------
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "base"
$ws.Range("C2").Value = 1577779312.350123
$ws.Range("D2").Value = "pathtopickle"
$ws.Range("E2").Value = "xlabel"
$ws.Range("F2").Value = "ylabel"

# Only A2 & B2 (the new two-level index columns) pick up the bold/bordered
# style; C2:F2 stay on the default style, matching the target sheet.
$ws.Range("A2:B2").PasteSpecial(-4122)

$excel.CutCopyMode = $false
